$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure D-column numeric-looking text values are written as TEXT (not auto-converted to numbers).
# We temporarily force Text number format on the whole D data range, write all the values,
# then restore the cells back to the default "Normal" style so the saved file keeps the original
# (un-styled) cell formatting, matching the source workbook.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "42.981.22"
$ws.Range("E2").Value = "  +0.12%  "
$ws.Range("D3").Value = "2.535.58"
$ws.Range("E3").Value = "  -1.09%  "
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.28%  "
$ws.Range("D5").Value = "305.60"
$ws.Range("E5").Value = "  +1.18%  "
$ws.Range("D6").Value = "100.87"
$ws.Range("E6").Value = "  +6.73%  "
$ws.Range("E7").Value = "  +1.02%  "
$ws.Range("E8").Value = "  +0.12%  "
$ws.Range("D9").Value = "0.546"
$ws.Range("E9").Value = "  +0.20%  "
$ws.Range("D10").Value = "37.62"
$ws.Range("E10").Value = "  +3.11%  "
$ws.Range("D11").Value = "0.0819"
$ws.Range("E11").Value = "  +1.04%  "
$ws.Range("D12").Value = "7.76"
$ws.Range("E12").Value = "  -0.06%  "
$ws.Range("D13").Value = "0.113"
$ws.Range("E13").Value = "  -0.41%  "
$ws.Range("D14").Value = "2.925.84"
$ws.Range("E14").Value = "  -0.95%  "
$ws.Range("D15").Value = "2.565.42"
$ws.Range("E15").Value = "  +0.27%  "
$ws.Range("D16").Value = "15.28"
$ws.Range("E16").Value = "  +7.59%  "
$ws.Range("D17").Value = "0.866"
$ws.Range("E17").Value = "  -1.72%  "
$ws.Range("D18").Value = "42.988.08"
$ws.Range("E18").Value = "  +0.04%  "
$ws.Range("D19").Value = "13.14"
$ws.Range("E19").Value = "  +3.30%  "
$ws.Range("D20").Value = "0.0₃0988"
$ws.Range("E20").Value = "  -1.10%  "
$ws.Range("E21").Value = "  -0.81%  "
$ws.Range("D22").Value = "71.70"
$ws.Range("E22").Value = "  +0.20%  "
$ws.Range("D23").Value = "254.20"
$ws.Range("E23").Value = "  +0.36%  "
$ws.Range("E24").Value = "  +0.07%  "
$ws.Range("E25").Value = "  -3.27%  "
$ws.Range("D26").Value = "27.28"
$ws.Range("E26").Value = "  -4.82%  "
$ws.Range("E27").Value = "  +0.17%  "
$ws.Range("D30").Value = "38.83"
$ws.Range("E30").Value = "  +4.67%  "
$ws.Range("D31").Value = "6.14"
$ws.Range("E31").Value = "  +1.31%  "
$ws.Range("D32").Value = "158.25"
$ws.Range("E32").Value = "  +2.96%  "
$ws.Range("E33").Value = "  -1.62%  "
$ws.Range("D34").Value = "0.0799"
$ws.Range("E34").Value = "  -0.14%  "
$ws.Range("E35").Value = "  -2.27%  "
$ws.Range("E36").Value = "  -3.60%  "
$ws.Range("D37").Value = "18.46"
$ws.Range("E37").Value = "  +2.97%  "
$ws.Range("E38").Value = "  +1.58%  "
$ws.Range("E39").Value = "  +0.10%  "
$ws.Range("D40").Value = "23.98"
$ws.Range("E40").Value = "  +3.58%  "
$ws.Range("D41").Value = "3.48"
$ws.Range("E41").Value = "  +1.74%  "
$ws.Range("D42").Value = "2.08"
$ws.Range("E42").Value = "  +1.06%  "
$ws.Range("D43").Value = "3.88"
$ws.Range("E43").Value = "  +0.24%  "
$ws.Range("E44").Value = "  -1.61%  "
$ws.Range("D45").Value = "0.999"
$ws.Range("E45").Value = "  +0.02%  "
$ws.Range("D46").Value = "2.046.10"
$ws.Range("E46").Value = "  -2.58%  "
$ws.Range("D47").Value = "86.19"
$ws.Range("E47").Value = "  +1.01%  "
$ws.Range("D48").Value = "9.00"
$ws.Range("E48").Value = "  -2.67%  "
$ws.Range("D49").Value = "2.783.31"
$ws.Range("E49").Value = "  -0.99%  "
$ws.Range("E50").Value = "  +1.41%  "
$ws.Range("D51").Value = "103.48"
$ws.Range("E51").Value = "  -2.94%  "

# Row 28/29: Toncoin and Cosmos swap rank order (new coin data for each).
$ws.Range("B28").Value = "Cosmos"
$ws.Range("C28").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D28").Value = "10.46"
$ws.Range("E28").Value = "  +1.98%  "

$ws.Range("B29").Value = "Toncoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D29").Value = "2.33"
$ws.Range("E29").Value = "  +9.27%  "

# Restore default styling on the D column (removes the temporary text-number-format).
$ws.Range("D2:D51").Style = "Normal"
